$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Sep 24, 2024"
$ws.Range("B5").Value = 74200
$ws.Range("C5").Value = 10487.93
$ws.Range("D5").Value = 9281.35
$ws.Range("E5").Value = 7.0168
